$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$ws.Range("G3").Value = 500
$ws.Range("G4").Value = 1300
$ws.Range("G8").Value = 2500
$ws.Range("G9").Value = 1500
$ws.Range("G10").Value = 900
$ws.Range("G12").Value = 1500
$ws.Range("G15").Value = 1000
$ws.Range("G20").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("G24").Value = 500
$ws.Range("G27").Value = 5000
$ws.Range("G30").Value = 0
$ws.Range("G31").Value = 1000
$ws.Range("G88").Value = 1900
$ws.Range("G89").Value = 15100
$ws.Range("G90").Value = 9000
$ws.Range("G91").Value = 2450
$ws.Range("G92").Value = 2450
$ws.Range("G93").Value = 750
$ws.Range("G94").Value = 0
$ws.Range("G95").Value = 1350
$ws.Range("G96").Value = 1300
$ws.Range("G98").Value = 650
$ws.Range("G99").Value = 5100
$ws.Range("G100").Value = 3935.11
$ws.Range("G101").Value = 7410
$ws.Range("G103").Value = 4330
$ws.Range("G105").Value = 950
$ws.Range("G106").Value = 0
$ws.Range("G107").Value = 1500
$ws.Range("G109").Value = 300
$ws.Range("G110").Value = 16110
$ws.Range("G111").Value = 300
$ws.Range("G114").Value = 5850
$ws.Range("G115").Value = 260
$ws.Range("G116").Value = 350
$ws.Range("G117").Value = 3100
$ws.Range("G119").Value = 0
$ws.Range("G121").Value = 18200
$ws.Range("G122").Value = 5600
$ws.Range("G123").Value = 100
$ws.Range("G126").Value = 5505
$ws.Range("G127").Value = 1020
$ws.Range("G128").Value = 2250
$ws.Range("G129").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("G132").Value = 3300
$ws.Range("G133").Value = 3750
$ws.Range("G134").Value = 300
$ws.Range("G135").Value = 550
$ws.Range("G138").Value = 1650
$ws.Range("G139").Value = 1000
$ws.Range("G140").Value = 0
$ws.Range("G144").Value = 3000
$ws.Range("G145").Value = 3000
$ws.Range("G146").Value = 1000
$ws.Range("G147").Value = 6500
$ws.Range("G159").Value = 1500
$ws.Range("G162").Value = 3000
$ws.Range("G163").Value = 1000
$ws.Range("G169").Value = 8000
$ws.Range("G171").Value = 3000
$ws.Range("G173").Value = 0
$ws.Range("G175").Value = 3000
$ws.Range("G178").Value = 1000
$ws.Range("G198").Value = 0
$ws.Range("G199").Value = 2000
$ws.Range("G203").Value = 2000
$ws.Range("G204").Value = 4500
$ws.Range("G208").Value = 2000
$ws.Range("G209").Value = 0
$ws.Range("G210").Value = 3500
$ws.Range("G212").Value = 0
$ws.Range("G214").Value = 1000
$ws.Range("G215").Value = 5000
$ws.Range("G218").Value = 2000
$ws.Range("G220").Value = 2500
$ws.Range("G227").Value = 1000
$ws.Range("G229").Value = 3000
$ws.Range("G232").Value = 3500
$ws.Range("G233").Value = 7500
$ws.Range("G234").Value = 5000
$ws.Range("G236").Value = 3000
$ws.Range("G263").Value = 1500
$ws.Range("G264").Value = 300
$ws.Range("G265").Value = 1000
$ws.Range("G266").Value = 8000
$ws.Range("G268").Value = 1150
$ws.Range("G271").Value = 1500
$ws.Range("G272").Value = 6000
$ws.Range("G273").Value = 2000
$ws.Range("G274").Value = 0
$ws.Range("G276").Value = 1200
$ws.Range("G279").Value = 1000
$ws.Range("G280").Value = 2000
$ws.Range("G284").Value = 400
$ws.Range("G285").Value = 2500
$ws.Range("G286").Value = 200
$ws.Range("G287").Value = 2890
$ws.Range("G288").Value = 1000
$ws.Range("G293").Value = 500
$ws.Range("G294").Value = 462561.11
